$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-24 Tuesday" "2024-09-25 Wednesday"

Replace-Text "596×9=5364" "218×3=654"
Replace-Text "680×5=3400" "763×2=1526"
Replace-Text "337×2=674" "208×7=1456"
Replace-Text "998×2=1996" "406×7=2842"
Replace-Text "619×4=2476" "777×9=6993"

Replace-Text "663×5=3315" "863×8=6904"
Replace-Text "332×4=1328" "462×4=1848"
Replace-Text "587×6=3522" "438×2=876"
Replace-Text "440×6=2640" "336×8=2688"
Replace-Text "870×7=6090" "311×9=2799"

Replace-Text "786×7=5502" "340×3=1020"
Replace-Text "536×7=3752" "719×7=5033"
Replace-Text "419×8=3352" "552×8=4416"
Replace-Text "143×5=715" "822×7=5754"
Replace-Text "343×9=3087" "342×4=1368"

Replace-Text "310×6=1860" "871×5=4355"
Replace-Text "683×2=1366" "464×2=928"
Replace-Text "507×5=2535" "432×9=3888"
Replace-Text "735×9=6615" "760×7=5320"
Replace-Text "198×6=1188" "550×9=4950"

Replace-Text "299×9=2691" "428×8=3424"
Replace-Text "225×6=1350" "630×5=3150"
Replace-Text "484×7=3388" "263×7=1841"
Replace-Text "607×7=4249" "320×3=960"
Replace-Text "405×7=2835" "831×5=4155"
